$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data per latest scrape
$ws.Range("D2").Value = "45.245.42"
$ws.Range("E2").Value = "  +2.42%  "
$ws.Range("D3").Value = "2.419.47"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.16"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.90"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.90%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.513"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  +4.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.38"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.17%  "
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("E12").Value = "  -1.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.22"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.74%  "
$ws.Range("E14").Value = "  +1.25%  "
$ws.Range("D15").Value = "2.800.04"
$ws.Range("E15").Value = "  -0.16%  "
$ws.Range("D16").Value = "2.420.05"
$ws.Range("E16").Value = "  -0.89%  "
$ws.Range("E17").Value = "  +1.01%  "
$ws.Range("D18").Value = "45.182.44"
$ws.Range("E18").Value = "  +2.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.21"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.32"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.22%  "
$ws.Range("D21").Value = "0.0₃0920"
$ws.Range("E21").Value = "  +1.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.72"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "244.26"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.60%  "
$ws.Range("E24").Value = "  -1.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.49"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.22%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.70"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.61"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "49.45"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.84"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.17"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +8.16%  "
$ws.Range("E33").Value = "  +7.33%  "
$ws.Range("E34").Value = "  +0.59%  "
$ws.Range("E35").Value = "  +0.32%  "
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("E37").Value = "  -1.31%  "
$ws.Range("E38").Value = "  -0.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "128.32"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.25%  "
$ws.Range("E40").Value = "  -0.60%  "
$ws.Range("E41").Value = "  +0.65%  "
$ws.Range("E42").Value = "  -4.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.55"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.24%  "
$ws.Range("E44").Value = "  +0.61%  "
$ws.Range("D45").Value = "1.940.41"
$ws.Range("E45").Value = "  -0.59%  "
$ws.Range("E46").Value = "  -2.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.93"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.14%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.79"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +9.02%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.12"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "76.83"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +4.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.80"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +5.15%  "

Write-Host "Applied 81 cell updates"
